$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de Usuario")
Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
